$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the formatting of H1
# (bold font, thin border, centered alignment) by copying H1's format over.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J (rows 2-12)
$data = @{
    2  = @(7, 7)
    3  = @(5, 5)
    4  = @(7, 7)
    5  = @(9, 9)
    6  = @(7, 7)
    7  = @(8, 8)
    8  = @(5, 6)
    9  = @(5, 5)
    10 = @(4, 4)
    11 = @(7, 7)
    12 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
